$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns hold numeric-looking
# text (e.g. "1.00", "65.599.75", "  -4.40%  "). Force those cells
# to Text format before assigning so Excel keeps the exact string
# instead of re-parsing it as a number (which would drop trailing
# zeros, switch to scientific notation, turn "+0.00%" into 0, etc).
$numericTextCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","E7","E8","D9","E9","E10","E11","D12","E12","E13","D14","E14","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","E22","D23","E23","D24","E24","E25","D26","E26","D27","E27","D28","E28","D29","E29","D30","E30","D31","E31","D32","E32","D33","E33","D34","E34","E35","D36","E36","E37","D38","E38","D39","E39","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","E48","E49","D50","E50","E51")
foreach ($cellRef in $numericTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values scraped on Tue Apr  2 17:04:18 UTC 2024
$ws.Range('D2').Value = '65.599.75'
$ws.Range('E2').Value = '  -4.40%  '
$ws.Range('D3').Value = '3.263.41'
$ws.Range('E3').Value = '  -5.35%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '555.24'
$ws.Range('E5').Value = '  -2.86%  '
$ws.Range('D6').Value = '180.70'
$ws.Range('E6').Value = '  -4.78%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  -2.83%  '
$ws.Range('D9').Value = '3.258.15'
$ws.Range('E9').Value = '  -5.22%  '
$ws.Range('E10').Value = '  -8.05%  '
$ws.Range('E11').Value = '  -4.43%  '
$ws.Range('D12').Value = '47.28'
$ws.Range('E12').Value = '  -7.01%  '
$ws.Range('E13').Value = '  -6.41%  '
$ws.Range('D14').Value = '635.22'
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('E15').Value = '  -5.35%  '
$ws.Range('D16').Value = '3.780.37'
$ws.Range('E16').Value = '  -5.61%  '
$ws.Range('D17').Value = '65.605.05'
$ws.Range('E17').Value = '  -4.68%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '0.116'
$ws.Range('E18').Value = '  -3.13%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '17.71'
$ws.Range('E19').Value = '  -1.83%  '
$ws.Range('D20').Value = '3.253.16'
$ws.Range('E20').Value = '  -5.88%  '
$ws.Range('D21').Value = '11.34'
$ws.Range('E21').Value = '  -7.33%  '
$ws.Range('E22').Value = '  -3.67%  '
$ws.Range('D23').Value = '17.74'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').Value = '105.93'
$ws.Range('E24').Value = '  +7.13%  '
$ws.Range('E25').Value = '  -6.77%  '
$ws.Range('D26').Value = '3.97'
$ws.Range('E26').Value = '  -6.66%  '
$ws.Range('D27').Value = '2.66'
$ws.Range('E27').Value = '  -5.56%  '
$ws.Range('D28').Value = '9.53'
$ws.Range('E28').Value = '  -2.47%  '
$ws.Range('D29').Value = '8.68'
$ws.Range('E29').Value = '  -4.89%  '
$ws.Range('D30').Value = '30.24'
$ws.Range('E30').Value = '  -5.95%  '
$ws.Range('D31').Value = '4.10'
$ws.Range('E31').Value = '  -0.94%  '
$ws.Range('D32').Value = '6.30'
$ws.Range('E32').Value = '  -5.55%  '
$ws.Range('D33').Value = '11.02'
$ws.Range('E33').Value = '  -4.23%  '
$ws.Range('D34').Value = '552.28'
$ws.Range('E34').Value = '  +9.94%  '
$ws.Range('E35').Value = '  -2.43%  '
$ws.Range('D36').Value = '57.05'
$ws.Range('E36').Value = '  -6.15%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').Value = '3.626.65'
$ws.Range('E38').Value = '  -0.71%  '
$ws.Range('D39').Value = '3.68'
$ws.Range('E39').Value = '  +6.35%  '
$ws.Range('E40').Value = '  -2.36%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.132'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = '2.74'
$ws.Range('E42').Value = '  -4.55%  '
$ws.Range('D43').Value = '0.0₃0716'
$ws.Range('E43').Value = '  -7.80%  '
$ws.Range('D44').Value = '31.90'
$ws.Range('E44').Value = '  -6.53%  '
$ws.Range('D45').Value = '0.337'
$ws.Range('E45').Value = '  -7.64%  '
$ws.Range('D46').Value = '3.32'
$ws.Range('E46').Value = '  -0.94%  '
$ws.Range('D47').Value = '0.0414'
$ws.Range('E47').Value = '  -4.40%  '
$ws.Range('E48').Value = '  -6.25%  '
$ws.Range('E49').Value = '  -3.41%  '
$ws.Range('D50').Value = '0.997'
$ws.Range('E50').Value = '  -0.42%  '
$ws.Range('E51').Value = '  +2.31%  '

Write-Output "Updated $($numericTextCells.Count + 8) cells."
